$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the last data row (row 96: "Create Country", "PASSED", "chrome", "07.04.23")
# into a new row 97, mirroring the repeated scenario row seen after the Cucumber/Jenkins
# example re-run merge. Copy/PasteSpecial preserves the original cell types (shared-string
# text, including the date-like "07.04.23" value) instead of letting it be reinterpreted.
$ws.Range("A96:D96").Copy()
$ws.Range("A97:D97").PasteSpecial()
$excel.CutCopyMode = $false
